# Auto update stock data
#
# Rolls the "as of" date (column A) forward one day, 2025/11/01 -> 2025/11/02,
# for every ticker's most-recent-data row, and refreshes the two EBITDA
# figures (column B) that were revised alongside that date roll.
#
# Cells in this sheet store their values as literal text (inline strings),
# even though some look like dates/numbers. Plain `.Value` assignment lets
# Excel auto-coerce a string like "2025/11/02" into a date serial (and
# "12.50" into the number 12.5), which would silently change the cell type
# and drop formatting/trailing zeros. To keep the values as plain text we
# temporarily force a text number format before assigning, then clear the
# format again afterwards (these cells had no explicit formatting to begin
# with) so the resulting style matches the original "General" default.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$dateRows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $dateRows) {
    Set-TextValue $ws.Cells.Item($r, 1) "2025/11/02"
}

# EBITDA corrections that came in with the refreshed date
Set-TextValue $ws.Cells.Item(20, 2) "12.50"
Set-TextValue $ws.Cells.Item(68, 2) "13.02"
